# Saldo_guide.xlsx update
# - Shift the "Dt. Referencia" (column G) date for every data row
#   from 2024-10-01 (serial 45566) to 2024-10-02 (serial 45567).
# - A handful of "Saldo Previsto"/"Vl. Total" (columns E & H) values were
#   refreshed for specific accounts (rows 109, 115, 148, 255).
# - The worksheet tab name encodes the report's generation timestamp; bump
#   it from 20241001-093545 to 20241002-085714 to match the new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 274 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq 45566) {
        $cell.Value2 = 45567
    }
}

# Targeted value refreshes (columns E = Saldo Previsto, H = Vl. Total)
$updates = @{
    109 = 41447.71
    115 = 24962.32
    148 = 0
    255 = 15841.4
}

foreach ($row in $updates.Keys) {
    $val = $updates[$row]
    $ws.Cells.Item($row, 5).Value2 = $val
    $ws.Cells.Item($row, 8).Value2 = $val
}

# Rename the sheet to reflect the newer export timestamp embedded in its name
$ws.Name = "IClientBalance-20241002-085714-"
